$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data values, preserving text cell type
# (these cells look numeric but are stored as text, so a leading apostrophe
# forces Excel to keep them as text instead of auto-converting to numbers).
$ws.Range("D2").Value = "'278.05"
$ws.Range("E2").Value = "'0.68%"
$ws.Range("D3").Value = "'27.24"
$ws.Range("E3").Value = "'1.94%"
$ws.Range("E4").Value = "'-0.14%"
$ws.Range("D5").Value = "'0.06427"
$ws.Range("E5").Value = "'1.55%"
$ws.Range("D6").Value = "'6.994"
$ws.Range("D7").Value = "'1.204"
$ws.Range("E7").Value = "'-8.88%"
$ws.Range("D8").Value = "'0.8861"
$ws.Range("E8").Value = "'1.23%"
$ws.Range("D9").Value = "'0.1550"
$ws.Range("E9").Value = "'-0.02%"
$ws.Range("D10").Value = "'0.05102"
$ws.Range("E10").Value = "'1.81%"
$ws.Range("D11").Value = "'0.07512"
$ws.Range("E11").Value = "'0.45%"
$ws.Range("D12").Value = "'0.02888"
$ws.Range("E12").Value = "'-1.48%"
$ws.Range("D13").Value = "'0.08973"
$ws.Range("E13").Value = "'-0.84%"
$ws.Range("D14").Value = "'0.001581"
$ws.Range("E14").Value = "'0.41%"
$ws.Range("D15").Value = "'0.0006364"
$ws.Range("E15").Value = "'0.64%"
$ws.Range("D16").Value = "'0.006128"
$ws.Range("E16").Value = "'1.59%"
$ws.Range("D17").Value = "'3.481"
$ws.Range("E17").Value = "'1.01%"
$ws.Range("D18").Value = "'3.305"
$ws.Range("E18").Value = "'-0.36%"
$ws.Range("D19").Value = "'2.218"
$ws.Range("E19").Value = "'-2.88%"
$ws.Range("E21").Value = "'1.24%"
$ws.Range("D22").Value = "'3.912"
$ws.Range("E22").Value = "'0.33%"
$ws.Range("D23").Value = "'0.04418"
$ws.Range("E23").Value = "'1.56%"
$ws.Range("E24").Value = "'8.70%"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("E25").Value = "'0.49%"
$ws.Range("D26").Value = "'0.003877"
$ws.Range("E26").Value = "'-7.96%"
$ws.Range("E28").Value = "'-1.61%"
$ws.Range("D29").Value = "'0.0001644"
$ws.Range("E29").Value = "'1.80%"
$ws.Range("D40").Value = "'0.04126"
$ws.Range("E40").Value = "'0.76%"
$ws.Range("D41").Value = "'0.006771"
$ws.Range("E41").Value = "'-2.86%"
$ws.Range("E42").Value = "'-0.01%"
$ws.Range("D43").Value = "'0.001901"
$ws.Range("E43").Value = "'-16.62%"
$ws.Range("D44").Value = "'0.01184"
$ws.Range("E44").Value = "'9.89%"
$ws.Range("D45").Value = "'0.00005330"
$ws.Range("E45").Value = "'0.72%"
$ws.Range("D46").Value = "'1.680"
$ws.Range("E46").Value = "'12.79%"
$ws.Range("D47").Value = "'0.01853"
$ws.Range("E47").Value = "'-7.34%"
